# FMEA workbook update: add "Safety Requirement" column (G) with
# safety-requirement text for each failure-mode row, per commit
# "Safety requirements for data updated".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Enter the new text values first, in the same order the original
# author typed them (G3, G2, G4, G1) so the shared-string table and
# row->value wiring line up with the authored workbook. ---
$ws.Range("G3").Value = "All data to be backed up by host software"
$ws.Range("G2").Value = "All users required to complete test prompt"
$ws.Range("G4").Value = "Only admin users will have access to device manager"
$ws.Range("G1").Value = "Safety Requirement"

# Merge the G4:G5 pair first (Device manager rows share one Safety
# Requirement cell) so the subsequent format copy lands on the final,
# merged pair of cells.
$ws.Range("G4:G5").Merge()

# --- Apply formatting by copying from the matching existing cells so the
# new column reuses the workbook's existing styles (header style, bordered
# left-aligned content style, bordered centered category style). ---

# G1 header -> same style as the other header cells (F1)
$ws.Range("F1").Copy()
$ws.Range("G1").PasteSpecial(-4122)

# G2:G3 content cells -> same style as other left-aligned bordered content cells (B2)
$ws.Range("B2").Copy()
$ws.Range("G2:G3").PasteSpecial(-4122)

# G4:G5 (merged) category-style cell -> same style as other centered bordered cells (A2)
$ws.Range("A2").Copy()
$ws.Range("G4").PasteSpecial(-4122)
$ws.Range("A2").Copy()
$ws.Range("G5").PasteSpecial(-4122)

$ws.Application.CutCopyMode = $false

# Row 3 grows to fit the new wrapped text in G3
$ws.Rows(3).RowHeight = 45

# New column width for the Safety Requirement column
$ws.Columns("G").ColumnWidth = 21.16667

# Selection moved (as recorded by the author's session) to B8
$ws.Range("B8").Select()
